$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 21) following the existing data pattern.
$row = 21

$ws.Cells.Item($row, 1).Value = 27      # Id
$ws.Cells.Item($row, 2).Value = 14      # ClientId
$ws.Cells.Item($row, 3).Value = 13      # ProcessedGroupId
$ws.Cells.Item($row, 4).Value = 2       # MachineNo
$ws.Cells.Item($row, 5).Value = 1       # VehicleNo
$ws.Cells.Item($row, 6).Value = 2018    # DtYear
$ws.Cells.Item($row, 7).Value = 4       # DtMonth
$ws.Cells.Item($row, 8).Value = "NULL"  # DtDay
$ws.Cells.Item($row, 9).Value = "NULL"  # Dt
$ws.Cells.Item($row, 10).Value = "NULL" # KmFrom
$ws.Cells.Item($row, 11).Value = 220660 # KmTo
# Column L (InsDt) left blank for this row.

# Update the active selection to reflect the next empty row, as in the source workbook.
$ws.Range("A22").Select()
